$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -17
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -6
$ws.Range("F12").Value = -3
$ws.Range("F16").Value = -4
$ws.Range("F18").Value = -3
$ws.Range("F19").Value = 1
